# Generate Report for Handback
# Updates the "Latest Handback DateTime" (column K) for the first handback
# file (ad1ecfcc-9247-45bd-88b9-6231b4fabed8) on both the "zh-cn" and
# "de-de" language sheets, reflecting newly generated handback reports.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("K2").Value = "2016-11-29 05:14:50"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("K2").Value = "2016-11-29 05:15:09"
